# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1, style index 1:
# bold font, thin border, centered) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values - constant across all 49 data rows (rows 2-50)
$ws.Range("AD2:AD50").Value = 71
$ws.Range("AE2:AE50").Value = 91
$ws.Range("AF2:AF50").Value = 0
